# Add a new "solidity" column before "thickness_max_chord_ratio" on the
# "geometry" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geometry")

# Insert a new column at AN, shifting AN:AQ -> AO:AR
$ws.Columns.Item("AN").Insert()

# Populate the new column's header and data value
$ws.Range("AN1").Value = "solidity"
$ws.Range("AN1").Style = $ws.Range("AO1").Style

$ws.Range("AN2").Value = "[1.42997704 1.70997375]"
